$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B4 text from "Prediction" to "Forecast"
$ws.Range("B4").Value = "Forecast"

# Remove now-unused trailing empty rows 10 and 11
$ws.Rows("10:11").Delete()
